$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The yearly report rolls forward one period: drop the oldest "1396/12"
# column (currently column E) by shifting everything one column to the
# left, then open a fresh column I for the new "1401/12" period.
$ws.Range("E1:E31").Delete(-4159)          # xlShiftToLeft
$ws.Range("I1:I31").Insert(-4161)          # xlShiftToRight
$ws.Columns("I").ColumnWidth = $ws.Columns("H").ColumnWidth

# Headers for the new period column
$ws.Range("I8").Value = "دوازده ماهه منتهی به 1401/12"
$ws.Range("I24").Value = "دوازده ماهه منتهی به 1401/12"

# Expense table (rows 10-20) - values for the newly opened "1401/12" column
$ws.Range("I10").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("I14").Value = 200259
$ws.Range("I15").Value = 0
$ws.Range("I16").Value = 226433
$ws.Range("I17").Value = 5184773
$ws.Range("I18").Value = 0
$ws.Range("I19").Value = 13413643
$ws.Range("I20").Value = 19025108

# Headcount table (rows 26-27)
$ws.Range("I26").Value = 1641
$ws.Range("I27").Value = 2117
